$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values (row 1): O1=14, P1=15
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Copy header style (bold, centered, bordered) from N1 to O1:P1
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

# Data rows: set O and P values for rows 2-67
$data = @(
    @(2, -0.2223177065780145, -0.2208349667846984),
    @(3, 0.2375351962286026, 0.2379846538284714),
    @(4, 0.2260171786819176, 0.2281294748961772),
    @(5, -0.04316101893056959, -0.04356786170631457),
    @(6, 0.194913862800399, 0.1949823890319755),
    @(7, -0.3748707848744286, -0.3732416199668248),
    @(8, -0.1438264973696326, -0.140434087662533),
    @(9, -0.3286813307253982, -0.3257935160846002),
    @(10, 0.4064714900680202, 0.4062414305362388),
    @(11, -0.1756242875344897, -0.1747059149733661),
    @(12, -0.05302997331425419, -0.05532338056412966),
    @(13, -0.02177798145155087, -0.02289633234179956),
    @(14, 0.08625681857296183, 0.08217584247616561),
    @(15, -0.02514523864141563, -0.0324035567947244),
    @(16, 0.3703242427303057, 0.3629288121562859),
    @(17, 0.5037311951121004, 0.4940887417911207),
    @(18, -0.111955041094718, -0.1171076175627849),
    @(19, 0.3572413774582678, 0.3523782313346648),
    @(20, 0.2836098781032829, 0.2747569997747542),
    @(21, 0.5605296518428371, 0.5534069929304815),
    @(22, 0.3199844789287077, 0.3119386640458501),
    @(23, -0.09680389022117211, -0.1037573427783416),
    @(24, 2.503146946857157, 2.374341945780941),
    @(25, 0.316029885414329, 0.3139337246386419),
    @(26, 0.2072208951722884, 0.2003653559908369),
    @(27, 0.07931807958825421, 0.0735870473083719),
    @(28, 0.821651183838776, 0.8183325641997663),
    @(29, 2.213922387034509, 2.118134169082537),
    @(30, 0.6944958502377289, 0.6887890958043077),
    @(31, -0.4218532290351512, -0.4224587858804125),
    @(32, 0.5811359413009466, 0.5770526675796561),
    @(33, 0.7703975105409337, 0.7684239484752134),
    @(34, -0.794868231502162, -0.7968638755882068),
    @(35, 0.8099744482102668, 0.810795360924124),
    @(36, 0.7666416579422214, 0.7689013399858302),
    @(37, 0.731914393734743, 0.7341129995180055),
    @(38, 0.6532251442004886, 0.6529032367204255),
    @(39, 0.6078398429931524, 0.609278849043142),
    @(40, 0.7746928014663623, 0.7756748598578272),
    @(41, 0.5604058759457682, 0.5621187350297192),
    @(42, 0.6374842897592634, 0.638518098205563),
    @(43, 0.6773935323672896, 0.6779883765976827),
    @(44, 0.6946927537993101, 0.6967351167497198),
    @(45, 0.6977143625818072, 0.7018493154729543),
    @(46, -1.235039754981151, -1.235488061208797),
    @(47, -0.9580365188461252, -0.9583295274417054),
    @(48, -0.8286462557775713, -0.8275672707147442),
    @(49, -0.6013438327822442, -0.6005724942852744),
    @(50, -0.03803848692888032, -0.03851800278800818),
    @(51, -0.814162298228931, -0.8128034318043048),
    @(52, -0.814162298228931, -0.8128034318043048),
    @(53, -1.081922677648855, -1.081418725237735),
    @(54, -0.1471888059183806, -0.1467365675696949),
    @(55, -0.9750342860138852, -0.9754987227464911),
    @(56, -0.8524039369883915, -0.8536268952318539),
    @(57, -0.9064072725796619, -0.9087961837413618),
    @(58, -1.054145063135318, -1.055794784747956),
    @(59, -0.7610149143850736, -0.7603236353503936),
    @(60, -0.4123648239502389, -0.4124435238430491),
    @(61, 0.390411563177137, 0.3905692590068945),
    @(62, -1.156189986759502, -1.15884085964446),
    @(63, -0.5606197347114803, -0.5569549653823174),
    @(64, -0.8350710662647094, -0.8339556325570525),
    @(65, -0.02620372393526893, -0.02610186147524414),
    @(66, -0.7397555265098047, -0.7427503439929491),
    @(67, -0.7049175008610312, -0.7094331090757803)
)

foreach ($item in $data) {
    $row = $item[0]
    $oVal = $item[1]
    $pVal = $item[2]
    $ws.Cells.Item($row, 15).Value = $oVal
    $ws.Cells.Item($row, 16).Value = $pVal
}

Write-Host "Done"